$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 651
$ws.Range("J2").Value = 6880
$ws.Range("K2").Value = 38
$ws.Range("L2").Value = 1859
$ws.Range("M2").Value = 117
$ws.Range("N2").Value = 1249
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 30
$ws.Range("Q2").Value = 15
$ws.Range("R2").Value = 86
$ws.Range("S2").Value = 777
$ws.Range("T2").Value = 1313
$ws.Range("U2").Value = 79
$ws.Range("V2").Value = 10847
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 11051
$ws.Range("Y2").Value = 13
$ws.Range("Z2").Value = 173
$ws.Range("AA2").Value = 82
